# Weekly update: insert a new price record at the top of the data block
# (row 9), pushing all existing rows (9-21) down by one (to 10-22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 9 - this shifts rows 9:21 down
# to 10:22 and copies formatting (incl. the date style on column D) from
# the row above, same as Excel's native "Insert Copied/Above Cells".
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new weekly record.
$ws.Range("A9").Value = 11
$ws.Range("B9").Value = "Vega Monumental Concepción"
$ws.Range("C9").Value = "Bíobío"
$ws.Range("D9").Value = 45272
$ws.Range("E9").Value = 8
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100101
$ws.Range("H9").Value = "Berries"
$ws.Range("I9").Value = 100101004
$ws.Range("J9").Value = "Frambuesa"
$ws.Range("K9").Value = "Sin especificar"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 100
$ws.Range("N9").Value = 10000
$ws.Range("O9").Value = 11000
$ws.Range("P9").Value = 10500
$ws.Range("Q9").Value = "$/bandeja 2 kilos"
$ws.Range("R9").Value = "Región de Ñuble"
$ws.Range("S9").Value = 5250
$ws.Range("T9").Value = 2
